# Apply the authored change to 03-iostream/00-iostream.pptx:
#   1. Reorder slides: move the "cerr / cout" picture slide (currently the
#      3rd slide in the deck) down to become the 8th slide.
#   2. Refresh the auto-updating "datetimeFigureOut" date placeholder text
#      (slide master + every slide layout) from 2021/8/22 to 2021/9/21.

$p = $ppt.ActivePresentation

# --- 1. Slide reorder -----------------------------------------------------
# Before: ... , slide3 (cerr/cout photo), slide4, slide5, slide6, slide7, slide8, ...
# After:  ... , slide4, slide5, slide6, slide7, slide8, slide3 (cerr/cout photo), ...
$moving = $p.Slides.Item(3)
$moving.MoveTo(8)

# --- 2. Refresh the date placeholder on the slide master and every layout -
$newDate = "2021/9/21"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shape = $master.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        if ($shape.TextFrame.TextRange.Text -eq "2021/8/22") {
            $shape.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shape = $layout.Shapes.Item($i)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -eq "2021/8/22") {
                $shape.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}
